# Applies the restructuring described by the diff: after the sheet is
# trimmed from 28 to 23 rows (the old rows 24-28 are removed), the block
# that used to live in rows 10-23 is re-laid-out: the long "Disciplina
# integradora..." / "Docentes responsaveis" paragraphs are replaced by the
# four lecturer-code rows interleaved with new "Programa resumido:" /
# "Short syllabus:" / "Programa:" / "Syllabus:" / "Avaliacao:" / "Metodo:" /
# "Criterio:" / "Norma de recuperacao:" / "Bibliografia:" / "Requisitos:"
# labels, and the two weak-requirement rows (LOQ4010 / LOQ4031) take the
# place of the old "Metodo:"/"Criterio:" evaluation-text rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 rows that no longer exist after the restructuring (old rows 24-28)
$ws.Range("24:28").Delete()

# Row 10
$ws.Range('A10').Value = 'Objetivos:'
$ws.Range('B10').Value = '198273 - Domingos Savio Giordani'
$ws.Range('C10').Value = '198273 - Domingos Savio Giordani'
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Range('A11').Value = 'Objectives:'
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Range('A12').Value = 'Programa resumido:'
$ws.Range('B12').Value = '5817045 - Elisângela de Jesus Cândido Moraes'
$ws.Range('C12').Value = '5817045 - Elisângela de Jesus Cândido Moraes'
$ws.Rows.Item(12).RowHeight = 60

# Row 13
$ws.Range('A13').Value = 'Short syllabus:'
$ws.Range("B13").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range('A14').Value = 'Programa:'
$ws.Range('B14').Value = '5817344 - Livia Melo Carneiro'
$ws.Range('C14').Value = '5817344 - Livia Melo Carneiro'
$ws.Rows.Item(14).RowHeight = 120

# Row 15
$ws.Range('A15').Value = 'Syllabus:'
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range('A16').Value = 'Avaliação:'
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Rows.Item(16).RowHeight = 15

# Row 17
$ws.Range('A17').Value = 'Método:'
$ws.Range('B17').Value = '787307 - Luis Fernando Figueiredo Faria'
$ws.Range('C17').Value = '787307 - Luis Fernando Figueiredo Faria'
$ws.Rows.Item(17).RowHeight = 60

# Row 18
$ws.Range('A18').Value = 'Critério:'
$ws.Range("B18").Value = @"
Avaliação de Projeto: 
-Apresentações orais (pré-projeto, relatório preliminar, relatório final, ampliação de escala);
-Trabalhos escritos (relatório preliminar e relatório final);
-Avaliação pelos pares.
"@
$ws.Range("C18").Value = @"
Avaliação de Projeto: 
-Apresentações orais (pré-projeto, relatório preliminar, relatório final, ampliação de escala);
-Trabalhos escritos (relatório preliminar e relatório final);
-Avaliação pelos pares.
"@
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range('A19').Value = 'Norma de recuperação:'
$ws.Range("B19").Value = @"
Média Final = Nota de Projeto 
Média final mínima de aprovação = 5,0
"@
$ws.Range("C19").Value = @"
Média Final = Nota de Projeto 
Média final mínima de aprovação = 5,0
"@
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range('A20').Value = 'Bibliografia:'
$ws.Range('B20').Value = '(Prova escrita + Média Final)/2         Nota Final Mínima para Aprovação= 5,0'
$ws.Range('C20').Value = '(Prova escrita + Média Final)/2         Nota Final Mínima para Aprovação= 5,0'
$ws.Rows.Item(20).RowHeight = 120

# Row 21
$ws.Range('A21').Value = 'Requisitos:'
$ws.Range("B21").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Rows.Item(21).RowHeight = 15

# Row 22
$ws.Range("A22").ClearContents()
$ws.Range("B22").Value = @"
LOQ4010 -  Introdução à  Engenharia  Química  (Requisito fraco)

"@
$ws.Range("C22").Value = @"
LOQ4010 -  Introdução à  Engenharia  Química  (Requisito fraco)

"@
$ws.Rows.Item(22).RowHeight = 30

# Row 23
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = @"
LOQ4031 -  Química Geral I  (Requisito fraco)

"@
$ws.Range("C23").Value = @"
LOQ4031 -  Química Geral I  (Requisito fraco)

"@
$ws.Rows.Item(23).RowHeight = 30

